$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'capri hex'
$ws.Cells.Item(2, 1).Value = 'knee pads softball'
$ws.Cells.Item(3, 1).Value = 'wrestling knee pads mens'
$ws.Cells.Item(4, 1).Value = 'snowboarding knee brace'
$ws.Cells.Item(5, 1).Value = 'basketball leggings for boys youth'
$ws.Cells.Item(6, 1).Value = 'honeycomb leggings'
$ws.Cells.Item(7, 1).Value = 'sliding knee pad'
$ws.Cells.Item(8, 1).Value = 'compression knee pad pants'
$ws.Cells.Item(9, 1).Value = 'youth basketball compression pants'
$ws.Cells.Item(10, 1).Value = 'basketball guide hand'
$ws.Cells.Item(11, 1).Value = '28 inch basketball'
$ws.Cells.Item(12, 1).Value = 'mens tights with knee pads'
$ws.Cells.Item(13, 1).Value = 'sliding sleeve softball'
$ws.Cells.Item(14, 1).Value = 'knee sleeve youth wrestling'
$ws.Cells.Item(15, 1).Value = 'basketball pads knee'
$ws.Cells.Item(16, 1).Value = 'sliding pants baseball youth'
$ws.Cells.Item(17, 1).Value = 'sports knee pads basketball'
$ws.Cells.Item(18, 1).Value = 'knee pads crossfit'
$ws.Cells.Item(19, 1).Value = 'boys knee pad leggings'
$ws.Cells.Item(20, 1).Value = 'knee pad basketball'
$ws.Cells.Item(21, 1).Value = 'basketball pads youth'
$ws.Cells.Item(22, 1).Value = 'knee pads for snowboarding'
$ws.Cells.Item(23, 1).Value = 'basketball leg sleeves for boys'
$ws.Cells.Item(24, 1).Value = 'basketball leggings youth boys'
$ws.Cells.Item(25, 1).Value = '20 30 compression leggings'
$ws.Cells.Item(26, 1).Value = 'workout knee pads'
$ws.Cells.Item(27, 1).Value = 'honeycomb compression pants'
$ws.Cells.Item(28, 1).Value = 'boys basketball leg sleeve'
$ws.Cells.Item(29, 1).Value = 'knee brace hockey'
$ws.Cells.Item(30, 1).Value = 'youth thigh compression sleeve'
$ws.Cells.Item(31, 1).Value = 'basketball tights youth boys'
$ws.Cells.Item(32, 1).Value = 'boys volleyball knee pads'
$ws.Cells.Item(33, 1).Value = 'crossfit pads'
$ws.Cells.Item(34, 1).Value = 'compression knee pads youth'
$ws.Cells.Item(35, 1).Value = 'padded basketball leggings for boys'
$ws.Cells.Item(36, 1).Value = 'crossfit knee sleeves men'
$ws.Cells.Item(37, 1).Value = 'baseball pants mens knee high'
$ws.Cells.Item(38, 1).Value = 'basketball knee pads for youth'
$ws.Cells.Item(39, 1).Value = 'softball knee pad'
$ws.Cells.Item(40, 1).Value = 'kneepads men'
$ws.Cells.Item(41, 1).Value = 'knee compression sleeve reduce strain'
$ws.Cells.Item(42, 1).Value = 'knee brace'
$ws.Cells.Item(43, 1).Value = 'knee compression sleeve'
$ws.Cells.Item(44, 1).Value = 'knee support'
$ws.Cells.Item(45, 1).Value = 'knee brace support'
$ws.Cells.Item(46, 1).Value = 'knee brace compression sleeve'
$ws.Cells.Item(47, 1).Value = 'compression sleeve'
$ws.Cells.Item(48, 1).Value = 'knee sleeve'
$ws.Cells.Item(49, 1).Value = 'knee compression'
$ws.Cells.Item(50, 1).Value = 'compression knee brace'
$ws.Cells.Item(51, 1).Value = 'knee pad'
$ws.Cells.Item(52, 1).Value = 'thigh compression sleeve'
$ws.Cells.Item(53, 1).Value = 'basketball training'
$ws.Cells.Item(54, 1).Value = 'crossfit training'
$ws.Cells.Item(55, 1).Value = 'basketball thigh pads'
$ws.Cells.Item(56, 1).Value = 'waist training leggings'
$ws.Cells.Item(57, 1).Value = 'work out tights mens'
$ws.Cells.Item(58, 1).Value = 'compression leggings capri'
$ws.Cells.Item(59, 1).Value = 'black workout leggings'
$ws.Cells.Item(60, 1).Value = 'workout legging'
$ws.Cells.Item(61, 1).Value = 'capri legging'
$ws.Cells.Item(62, 1).Value = 'workout tights men'
$ws.Cells.Item(63, 1).Value = 'workout pad'
$ws.Cells.Item(64, 1).Value = 'leggings basketball'
$ws.Cells.Item(65, 1).Value = 'black legging'
$ws.Cells.Item(66, 1).Value = 'basketball knee sleeve men'
$ws.Cells.Item(67, 1).Value = 'basketball pads for training'
$ws.Cells.Item(68, 1).Value = 'workout tights for men'
$ws.Cells.Item(69, 1).Value = 'knee pads basketball mcdavid black'
$ws.Cells.Item(70, 1).Value = 'basketball knee pads adult black'
$ws.Cells.Item(71, 1).Value = 'mens leggings compression'
$ws.Cells.Item(72, 1).Value = 'black capri legging'
$ws.Cells.Item(73, 1).Value = 'man leggings'
$ws.Cells.Item(74, 1).Value = 'basketball sleeve black'
$ws.Cells.Item(75, 1).Value = 'men workout pants'
$ws.Cells.Item(76, 1).Value = 'basketball knee pad'
$ws.Cells.Item(77, 1).Value = 'hip workout'
$ws.Cells.Item(78, 1).Value = 'compression legging'
$ws.Cells.Item(79, 1).Value = 'workout pant'
$ws.Cells.Item(80, 1).Value = 'capri pant'
$ws.Cells.Item(81, 1).Value = 'youth training basketball'
$ws.Cells.Item(82, 1).Value = 'compression pant'
$ws.Cells.Item(83, 1).Value = 'baseball pants knee high mens'
$ws.Cells.Item(84, 1).Value = 'basketball knee pads youth girls'
$ws.Cells.Item(85, 1).Value = 'basketball knee protector'
$ws.Cells.Item(86, 1).Value = 'basketball leggings with knee pads'
$ws.Cells.Item(87, 1).Value = 'basketball youth knee pads'
$ws.Cells.Item(88, 1).Value = 'compression knee pants'
$ws.Cells.Item(89, 1).Value = 'compression leggings boys basketball'
$ws.Cells.Item(90, 1).Value = 'knee guards for basketball'
$ws.Cells.Item(91, 1).Value = 'knee pad basketball men'
$ws.Cells.Item(92, 1).Value = 'knee pad pants basketball'
$ws.Cells.Item(93, 1).Value = 'knee pads basketball'
$ws.Cells.Item(94, 1).Value = 'knee pads boys basketball'
$ws.Cells.Item(95, 1).Value = 'knee pads compression pants'
$ws.Cells.Item(96, 1).Value = 'knee pads for basketball men'
$ws.Cells.Item(97, 1).Value = 'knee pads men basketball'
$ws.Cells.Item(98, 1).Value = 'knee pads pants for men'
$ws.Cells.Item(99, 1).Value = 'knee pads wrestling'
$ws.Cells.Item(100, 1).Value = 'knee protection for workout'
